$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data (and row content shifts for rows 15-25)
# Each target cell is forced to Text format ("@") before assignment so that
# numeric-looking strings (e.g. "307.87", "-4.30%") are preserved as literal
# text, matching the original inlineStr cell type used throughout the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '307.87'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.30%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.13'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-8.77%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.084'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.43%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07693'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-6.67%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.252'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.79%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.633'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-8.38%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.48%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1026'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-8.58%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1769'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-6.01%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09246'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.24%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04441'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.12%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1057'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.12%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001237'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-4.48%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005866'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '3.01%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.007491'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2,413.11%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.361'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.14%'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.433'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.52%'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3311'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.65%'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.844'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-8.11%'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1347'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-3.04%'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2814'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '10.32%'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04144'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.72%'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001201'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-4.10%'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'HotbitToken'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004099'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-4.36%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001297'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '6.25%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02456'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-7.76%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05180'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-7.95%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007930'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.78%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-5.95%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007150'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '9.41%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001946'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-8.15%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007967'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.87%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3073'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-11.92%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.01%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003002'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-26.84%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004486'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '33.25%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.01%'
